# "Created apps types test"
# Adds a new worksheet "SystemStructureCreateAppsTypes" at the end of the
# workbook, populated with a small header + 5 rows of test data, mirroring
# the existing "SystemStructureCreateOrgsApps"-style sheets.

$wb = $excel.ActiveWorkbook

# New sheet goes after the last existing sheet (becomes the new last/active tab).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "SystemStructureCreateAppsTypes"

# ---- Header row ----
$ws.Range("A1").Value = "#"
$ws.Range("B1").Value = "Наименование"
$ws.Range("C1").Value = "Идентификатор"
$ws.Range("D1").Value = "Является корпоративным"
$ws.Range("E1").Value = "Поддерживает работу с компонентами"
$ws.Range("A1:E1").Font.Bold = $true

# ---- Data rows (filled column by column, matching source authoring order) ----
$ws.Range("B2").Value = "Тестовое наименование 1"
$ws.Range("B3").Value = "Тестовое наименование 2"
$ws.Range("B4").Value = "Тестовое наименование 3"
$ws.Range("B5").Value = "Тестовое наименование 4"
$ws.Range("B6").Value = "Тестовое наименование 5"

$ws.Range("C2").Value = "Тестовый идентификатор 1"
$ws.Range("C3").Value = "Тестовый идентификатор 2"
$ws.Range("C4").Value = "Тестовый идентификатор 3"
$ws.Range("C5").Value = "Тестовый идентификатор 4"
$ws.Range("C6").Value = "Тестовый идентификатор 5"

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

$ws.Range("D2").Value = "'True"
$ws.Range("D3").Value = "'False"
$ws.Range("D4").Value = "'True"
$ws.Range("D5").Value = "'False"
$ws.Range("D6").Value = "'True"

$ws.Range("E2").Value = "'True"
$ws.Range("E3").Value = "'False"
$ws.Range("E4").Value = "'False"
$ws.Range("E5").Value = "'True"
$ws.Range("E6").Value = "'False"

# ---- Column widths (characters), matching the source sheet layout ----
$ws.Columns("B").ColumnWidth = 29.09
$ws.Columns("C").ColumnWidth = 26.59
$ws.Columns("D").ColumnWidth = 25.42
$ws.Columns("E").ColumnWidth = 37.59

# ---- View state ----
$ws.Range("B14").Select() | Out-Null
